# Populate Sheet1 with the thrust-stand test data (A1:C4).
# Column A holds zero-padded numeric codes, so it must be forced to Text
# format before assignment, otherwise Excel would strip the leading zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:A4").NumberFormat = "@"

$data = @(
    @("005005", 57, 57),
    @("000000", 57, 57),
    @("045005", 57, 57),
    @("000000", 57, 57)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

Write-Host "Populated A1:C4 with test data"
